{"js": "// The document carries a single \"_GoBack\" bookmark (Word's automatic\n// \"last edit location\" marker). It currently sits at the very end of the\n// \"EOS ID: 1109\" paragraph (after the \"Solution: Solution-placeholder 1109\"\n// run). This edit relocates it into the \"EOS ID: 1123\" paragraph, right\n// after the \"Problem: Problem-placeholder 1123\" run and before the line\n// break that precedes \"Solution: Solution-placeholder 1123\" \u2014 a basic\n// integrity-check marker left by the author's last edit in that spot.\n\n// 1) Remove the bookmark from its current location.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the anchor text that the bookmark should now follow.\nconst results = context.document.body.search(\"Problem: Problem-placeholder 1123\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Anchor text 'Problem: Problem-placeholder 1123' not found\");\n}\n\n// 3) Re-insert the bookmark immediately after that run, before the\n//    subsequent line break + \"Solution: ...\" run.\nconst targetRange = results.items[0].getRange(\"After\");\ntargetRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The document carries a single \"_GoBack\" bookmark (Word's automatic\n# \"last edit location\" marker). It currently sits at the very end of the\n# \"EOS ID: 1109\" paragraph (right after the \"Solution: Solution-placeholder\n# 1109\" run). This edit relocates it into the \"EOS ID: 1123\" paragraph,\n# right after the \"Problem: Problem-placeholder 1123\" run and before the\n# line break that precedes \"Solution: Solution-placeholder 1123\" - a basic\n# integrity-check marker for the author's last edit in that spot.\n\n$d = $word.ActiveDocument\n\n# Locate the end of the \"Problem: Problem-placeholder 1123\" run.\n$r = $d.Content\n$found = $r.Find.Execute(\"Problem: Problem-placeholder 1123\")\nif (-not $found) {\n    throw \"Anchor text 'Problem: Problem-placeholder 1123' not found\"\n}\n\n# Collapse the found range to its end point (wdCollapseEnd = 0) so the\n# bookmark is inserted right there, not spanning the found text.\n$r.Collapse(0)\n\n# Re-adding a bookmark with a name that already exists moves it to the\n# new range instead of creating a duplicate, so this both removes the\n# bookmark from its old spot (end of the 1109 paragraph) and places it\n# in the new one.\n$d.Bookmarks.Add(\"_GoBack\", $r)\n"}
